$d = $word.ActiveDocument

function New-PkgXml([string]$bodyXml) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $bodyXml + '</w:body></w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'
    return $xml
}

function Set-ParagraphRuns([int]$paraIndex, [string]$runsXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $pStart = $p.Range.Start
    $pEnd = $p.Range.End - 1
    $r = $d.Range($pStart, $pEnd)
    $pXml = "<w:p>" + $runsXml + "</w:p>"
    $frag = New-PkgXml $pXml
    $r.InsertXML($frag)
}

# ---------------------------------------------------------------------------
# Paragraph 1: Title
# ---------------------------------------------------------------------------
$rPrTitle = '<w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="44"/></w:rPr>'
$titleRuns = "<w:r>$rPrTitle<w:t>The Spine: A Journey Through History</w:t></w:r>"
Set-ParagraphRuns 1 $titleRuns

# ---------------------------------------------------------------------------
# Paragraph 2: Author name ("Dr." + "Benjamin Miller")
# ---------------------------------------------------------------------------
$rPrAuthor = '<w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr>'
$authorRuns = "<w:r>$rPrAuthor<w:t>Dr</w:t></w:r>" + `
    "<w:r>$rPrAuthor<w:t>.</w:t></w:r>" + `
    "<w:r>$rPrAuthor<w:t>Benjamin Miller</w:t></w:r>"
Set-ParagraphRuns 2 $authorRuns

# ---------------------------------------------------------------------------
# Paragraph 3: Email ("bmilleremail" + "." + "org")
# ---------------------------------------------------------------------------
$rPrEmail = '<w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr>'
$emailRuns = "<w:r>$rPrEmail<w:t>bmilleremail</w:t></w:r>" + `
    "<w:r>$rPrEmail<w:t>.</w:t></w:r>" + `
    "<w:r>$rPrEmail<w:t>org</w:t></w:r>"
Set-ParagraphRuns 3 $emailRuns

Write-Host "Stage 1 complete"

# ---------------------------------------------------------------------------
# Paragraph 5: main body text
# ---------------------------------------------------------------------------
$rPrBody = '<w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr>'

$bodyRuns = ""
$bodyRuns += "<w:r>$rPrBody<w:t>In an era dominated by technology, few truly understand the intricate web of bones and muscles that define our existence</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:t>.</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:t xml:space=""preserve""> Amidst the digital tapestry of our lives, it is crucial to not lose sight of our physical foundation</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:t>.</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:t xml:space=""preserve""> Embark upon an expedition into the past, uncovering the captivating journey of the spine, an enigmatic structure that has inspired unravelling its secrets</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:t>.</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:t xml:space=""preserve""> Join me as we traverse epochs, tracing the transformation of the spine from its earliest origins to its modern-day significance, uncovering the mysteries hidden within its robust framework</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:t>.</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:br/></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:br/><w:t>From ancient civilizations to the boundaries of modern medicine, the spine has been a constant enigma, a symphony of bone, cartilage, and nerves, intricate yet robust, yielding to the complexities of its function</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:t>.</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:t xml:space=""preserve""> Its narrative is not solely confined to the annals of scientific discovery, but rather a tapestry woven with culture, tradition, and artistry, spanning across civilizations</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:t>.</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:t xml:space=""preserve""> As our voyage unfolds, we shall uncover how the spine has bewitched poets, captivated artists, and captivated scientists alike</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:t>.</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:br/></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:br/><w:t>The spine, a central pillar of our anatomy, an architectural feat of evolution, serves as a protective sheath for the spinal cord, granting mobility and supporting the body's weight</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:t>.</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:t xml:space=""preserve""> As we embark on this intellectual odyssey, we shall delve into the realm of biology and uncover the marvels of the spine's structure, the intricate dance of nerves, muscles, and bones</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:t>.</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:t xml:space=""preserve""> Along this path, we shall endeavour to unravel the enigmatic relationship between the spine and human consciousness, exploring the marvels of perception, sensation, and motor control</w:t></w:r>"
$bodyRuns += "<w:r>$rPrBody<w:t>.</w:t></w:r>"

Set-ParagraphRuns 5 $bodyRuns

Write-Host "Stage 2 complete"
Write-Host $d.Paragraphs.Item(5).Range.Text
